$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$v = $ws.Range("E1").Value2
Write-Output $v
